$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (naive forecaster bug: stray y_0_forecast/y_1_forecast
# values written for rows that should not have had a forecast yet).
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Recomputed naive forecaster values (tiny floating point precision fix).
$ws.Range("E3").Value = 6.253707197847591
$ws.Range("E4").Value = -1.563633964192079
$ws.Range("C5").Value = -7.921319741078636
$ws.Range("E5").Value = -14.35806537048446
$ws.Range("C6").Value = -4.774178217057756
$ws.Range("E9").Value = 6.136355062499965
$ws.Range("E10").Value = 2.114249845651872
$ws.Range("C11").Value = 1.477633171193116
$ws.Range("C12").Value = 1.239479831392831
$ws.Range("C13").Value = -0.2005250704869121
$ws.Range("C14").Value = 0.2379616621360992
$ws.Range("C15").Value = 2.234527904461148
$ws.Range("C17").Value = 1.311727872618218
$ws.Range("E17").Value = 1.093673275363694
$ws.Range("C18").Value = 1.470039379455734
$ws.Range("E18").Value = 1.339087911421122
$ws.Range("C19").Value = 1.784808447869191
$ws.Range("C20").Value = 1.638797242243228
$ws.Range("E20").Value = 1.006353890555212
$ws.Range("C23").Value = 2.634902838428399
$ws.Range("E23").Value = 2.450082126686826
$ws.Range("C24").Value = 2.159589514946725
$ws.Range("C25").Value = 2.330656125352215
$ws.Range("E26").Value = -0.301339632123987
$ws.Range("C27").Value = 0.2120367165967307
$ws.Range("E27").Value = 0.06930121167905146
$ws.Range("C28").Value = 0.8014493436638848
$ws.Range("E28").Value = 1.609625625599986
$ws.Range("C29").Value = 0.5878492443567529
$ws.Range("E29").Value = -0.184185216762256
$ws.Range("E30").Value = 0.2691345740890139
$ws.Range("E31").Value = 0.1120143486733172
$ws.Range("C32").Value = -3.107661574595766
$ws.Range("E32").Value = -8.513835774400015
$ws.Range("C33").Value = -8.268943763593073
$ws.Range("E33").Value = -31.94180729997805
$ws.Range("E34").Value = 23.52713729381606
$ws.Range("C35").Value = 1.269625353117143
$ws.Range("C36").Value = -1.621578487659103
$ws.Range("C37").Value = 0.7868572467511825
$ws.Range("E37").Value = 5.737052250893782
$ws.Range("C38").Value = 1.099928004397577
$ws.Range("E38").Value = 6.182044950645027
$ws.Range("C40").Value = 1.906376895025041
$ws.Range("E40").Value = 0.9449384537270955
$ws.Range("C42").Value = 2.310042359896247
$ws.Range("C43").Value = -0.4844381132618314
$ws.Range("E43").Value = -1.729907735624059
$ws.Range("C44").Value = -0.4925007786849234
$ws.Range("C45").Value = -0.0462759835394233
$ws.Range("E45").Value = 0.04644557489565937
$ws.Range("E46").Value = -0.3513551123189074
$ws.Range("C47").Value = -0.6828258883444516
$ws.Range("E47").Value = -1.143293480177665
$ws.Range("C49").Value = -0.243463903689245
$ws.Range("E49").Value = -0.1929563298375014
$ws.Range("C50").Value = -0.3101476031197037
$ws.Range("E50").Value = 0.2561130241983456
$ws.Range("C51").Value = -0.5104822884906102
$ws.Range("E51").Value = -0.799743249960494
$ws.Range("E53").Value = -0.9576579957581766
